$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-30 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-31 Sunday", 2) | Out-Null
$d.Content.Find.Execute("270×5=", $true, $false, $false, $false, $false, $true, 1, $false, "382×9=", 2) | Out-Null
$d.Content.Find.Execute("749×4=", $true, $false, $false, $false, $false, $true, 1, $false, "452×4=", 2) | Out-Null
$d.Content.Find.Execute("824×4=", $true, $false, $false, $false, $false, $true, 1, $false, "529×6=", 2) | Out-Null
$d.Content.Find.Execute("986×7=", $true, $false, $false, $false, $false, $true, 1, $false, "524×9=", 2) | Out-Null
$d.Content.Find.Execute("883×5=", $true, $false, $false, $false, $false, $true, 1, $false, "947×6=", 2) | Out-Null
$d.Content.Find.Execute("625×2=", $true, $false, $false, $false, $false, $true, 1, $false, "437×4=", 2) | Out-Null
$d.Content.Find.Execute("596×2=", $true, $false, $false, $false, $false, $true, 1, $false, "909×7=", 2) | Out-Null
$d.Content.Find.Execute("780×3=", $true, $false, $false, $false, $false, $true, 1, $false, "908×7=", 2) | Out-Null
$d.Content.Find.Execute("890×7=", $true, $false, $false, $false, $false, $true, 1, $false, "736×5=", 2) | Out-Null
$d.Content.Find.Execute("980×3=", $true, $false, $false, $false, $false, $true, 1, $false, "158×4=", 2) | Out-Null
$d.Content.Find.Execute("938×3=", $true, $false, $false, $false, $false, $true, 1, $false, "617×2=", 2) | Out-Null
$d.Content.Find.Execute("189×9=", $true, $false, $false, $false, $false, $true, 1, $false, "993×8=", 2) | Out-Null
$d.Content.Find.Execute("315×4=", $true, $false, $false, $false, $false, $true, 1, $false, "152×9=", 2) | Out-Null
$d.Content.Find.Execute("125×7=", $true, $false, $false, $false, $false, $true, 1, $false, "664×5=", 2) | Out-Null
$d.Content.Find.Execute("431×3=", $true, $false, $false, $false, $false, $true, 1, $false, "551×9=", 2) | Out-Null
$d.Content.Find.Execute("966×5=", $true, $false, $false, $false, $false, $true, 1, $false, "221×2=", 2) | Out-Null
$d.Content.Find.Execute("255×6=", $true, $false, $false, $false, $false, $true, 1, $false, "270×2=", 2) | Out-Null
$d.Content.Find.Execute("536×5=", $true, $false, $false, $false, $false, $true, 1, $false, "329×4=", 2) | Out-Null
$d.Content.Find.Execute("939×5=", $true, $false, $false, $false, $false, $true, 1, $false, "307×4=", 2) | Out-Null
$d.Content.Find.Execute("413×7=", $true, $false, $false, $false, $false, $true, 1, $false, "683×9=", 2) | Out-Null
$d.Content.Find.Execute("897×3=", $true, $false, $false, $false, $false, $true, 1, $false, "838×6=", 2) | Out-Null
$d.Content.Find.Execute("149×4=", $true, $false, $false, $false, $false, $true, 1, $false, "606×5=", 2) | Out-Null
$d.Content.Find.Execute("912×9=", $true, $false, $false, $false, $false, $true, 1, $false, "735×8=", 2) | Out-Null
$d.Content.Find.Execute("889×2=", $true, $false, $false, $false, $false, $true, 1, $false, "611×9=", 2) | Out-Null
$d.Content.Find.Execute("647×2=", $true, $false, $false, $false, $false, $true, 1, $false, "612×9=", 2) | Out-Null
